{"js": "// The document contains a single table whose data rows (0, 4, 8, 12, 16 -\n// interleaved with blank rows) hold short \"NN\u00f7N=\" division exercises.\n// The edit updates the values of those cells. Row 0 keeps the same number\n// of populated cells (5) even though the underlying OOXML diff shows two\n// cells removed and two different cells inserted - net effect is simply\n// that the 5 visible values in that row change. We therefore perform the\n// edit as a set of per-cell text replacements, which preserves all\n// existing run/paragraph formatting (font, size, justification).\n\nconst rowEdits = {\n  0: [\"32\u00f73=\", \"34\u00f77=\", \"68\u00f77=\", \"66\u00f76=\", \"88\u00f72=\"],\n  4: [\"57\u00f72=\", \"17\u00f73=\", \"86\u00f72=\", \"15\u00f78=\", \"14\u00f78=\"],\n  8: [\"43\u00f72=\", \"51\u00f79=\", \"63\u00f78=\", \"66\u00f79=\", \"98\u00f76=\"],\n  12: [\"53\u00f73=\", \"38\u00f72=\", \"34\u00f74=\", \"97\u00f72=\", \"25\u00f72=\"],\n  16: [\"96\u00f76=\", \"88\u00f75=\", \"64\u00f77=\", \"19\u00f73=\", \"81\u00f77=\"],\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Collect the search-result ranges for every cell we need to change.\nconst targets = []; // { range, newText }\nfor (const rowIndexStr of Object.keys(rowEdits)) {\n  const rowIndex = Number(rowIndexStr);\n  const newValues = rowEdits[rowIndex];\n  const row = rows.items[rowIndex];\n\n  for (let c = 0; c < row.cells.items.length; c++) {\n    const cell = row.cells.items[c];\n    cell.body.load(\"text\");\n    targets.push({ cell, newText: newValues[c] });\n  }\n}\nawait context.sync();\n\nfor (const t of targets) {\n  const oldText = t.cell.body.text;\n  const results = t.cell.body.search(oldText, {\n    matchCase: true,\n    matchWildcards: false,\n  });\n  results.load(\"items\");\n  t.results = results;\n}\nawait context.sync();\n\nfor (const t of targets) {\n  if (t.results.items.length > 0) {\n    t.results.items[0].insertText(t.newText, \"Replace\");\n  } else {\n    // Fallback: cell had no existing run to anchor to - just insert text.\n    t.cell.body.insertText(t.newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# The document contains a single table whose data rows (Word COM 1-based\n# row numbers 1, 5, 9, 13, 17 - interleaved with blank rows) hold short\n# \"NN\u00f7N=\" division exercises. The edit updates the values of those cells.\n# Row 1 keeps the same number of populated cells (5) even though the\n# underlying OOXML diff shows two cells removed and two different cells\n# inserted - net effect is simply that the 5 visible values in that row\n# change. We therefore perform the edit as a set of per-cell text\n# assignments, which preserves existing run/paragraph formatting\n# (font, size, justification) because Word keeps the cell's trailing\n# paragraph mark (and its formatting) when only the Range.Text is set.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$rowEdits = @{\n    1  = @(\"32\u00f73=\", \"34\u00f77=\", \"68\u00f77=\", \"66\u00f76=\", \"88\u00f72=\")\n    5  = @(\"57\u00f72=\", \"17\u00f73=\", \"86\u00f72=\", \"15\u00f78=\", \"14\u00f78=\")\n    9  = @(\"43\u00f72=\", \"51\u00f79=\", \"63\u00f78=\", \"66\u00f79=\", \"98\u00f76=\")\n    13 = @(\"53\u00f73=\", \"38\u00f72=\", \"34\u00f74=\", \"97\u00f72=\", \"25\u00f72=\")\n    17 = @(\"96\u00f76=\", \"88\u00f75=\", \"64\u00f77=\", \"19\u00f73=\", \"81\u00f77=\")\n}\n\nforeach ($rowIndex in $rowEdits.Keys) {\n    $values = $rowEdits[$rowIndex]\n    for ($c = 1; $c -le $values.Count; $c++) {\n        $cell = $tbl.Cell($rowIndex, $c)\n        $cell.Range.Text = $values[$c - 1]\n    }\n}\n"}
